# "Reverse Linked List easy" commit:
#  - DW sheet, row 23's problem name is trimmed from
#    "Distinct Palindrom Substring" to " Palindrom Substring"
#  - a brand-new row (32) is appended for a GeeksForGeeks variant of the
#    same problem ("Distinct Palindrom Substring" again, new link/remark)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DW")

# --- Row 23: rename the problem statement cell ---------------------------
$ws.Range("A23").Value = " Palindrom Substring"

# --- Row 32: new problem entry --------------------------------------------
$ws.Range("A32").Value = "Distinct Palindrom Substring"
$ws.Range("B32").Value = "String"
$ws.Range("C32").Value = "Medium"
$ws.Range("D32").Value = "https://practice.geeksforgeeks.org/problems/distinct-palindromic-substrings"
$ws.Range("E32").Value = "Just add a list to check for repeating "

# Row 22 carries the same "Neutral"-style banding (yellow fill / brown text,
# thin vertical borders on A:C, no border on D, bottom border on E) that the
# new row should have, so clone its per-cell formatting onto row 32.
$ws.Range("A22").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("B22").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("C22").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("D22").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E22").Copy()
$ws.Range("E32").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match the saved selection: the whole new row selected.
$ws.Range("A32:XFD32").Select()
